# Update the multiplication problems/answers in the single table.
# Cells are addressed directly by (row, column) rather than via
# Find/Replace because some new values coincide with other old values
# elsewhere in the table (e.g. "800x3=2400" and "746x9=6714" are both a
# replacement target in one spot and a replacement source elsewhere), so a
# global text replace could clobber a cell that was already updated.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "605×4=2420"
$t.Cell(1, 2).Range.Text = "815×2=1630"
$t.Cell(1, 3).Range.Text = "905×9=8145"
$t.Cell(1, 4).Range.Text = "931×8=7448"
$t.Cell(1, 5).Range.Text = "715×9=6435"

$t.Cell(5, 1).Range.Text = "257×5=1285"
$t.Cell(5, 2).Range.Text = "584×6=3504"
$t.Cell(5, 3).Range.Text = "746×9=6714"
$t.Cell(5, 4).Range.Text = "188×8=1504"
$t.Cell(5, 5).Range.Text = "106×8=848"

$t.Cell(10, 1).Range.Text = "514×2=1028"
$t.Cell(10, 2).Range.Text = "712×7=4984"
$t.Cell(10, 3).Range.Text = "716×6=4296"
$t.Cell(10, 4).Range.Text = "264×8=2112"
$t.Cell(10, 5).Range.Text = "571×3=1713"

$t.Cell(15, 1).Range.Text = "961×7=6727"
$t.Cell(15, 2).Range.Text = "800×3=2400"
$t.Cell(15, 3).Range.Text = "566×2=1132"
$t.Cell(15, 4).Range.Text = "372×4=1488"
$t.Cell(15, 5).Range.Text = "659×2=1318"

$t.Cell(20, 1).Range.Text = "332×7=2324"
$t.Cell(20, 2).Range.Text = "858×2=1716"
$t.Cell(20, 3).Range.Text = "927×5=4635"
$t.Cell(20, 4).Range.Text = "956×7=6692"
$t.Cell(20, 5).Range.Text = "151×7=1057"
